$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.812.68"
$ws.Range("E2").Value = "  -1.92%  "

$ws.Range("D3").Value = "3.492.37"
$ws.Range("E3").Value = "  -1.34%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").Value = "'601.61"
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("D6").Value = "'148.06"
$ws.Range("E6").Value = "  -3.47%  "

$ws.Range("D7").Value = "3.490.64"
$ws.Range("E7").Value = "  -1.36%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.479"
$ws.Range("E9").Value = "  -2.26%  "

$ws.Range("E10").Value = "  -1.42%  "

$ws.Range("D11").Value = "'7.95"
$ws.Range("E11").Value = "  +3.97%  "

$ws.Range("E12").Value = "  -3.10%  "

$ws.Range("E13").Value = "  -2.60%  "

$ws.Range("D14").Value = "4.081.94"
$ws.Range("E14").Value = "  -1.36%  "

$ws.Range("D15").Value = "'31.24"
$ws.Range("E15").Value = "  -5.36%  "

$ws.Range("D16").Value = "3.490.84"
$ws.Range("E16").Value = "  -1.35%  "

$ws.Range("D17").Value = "66.778.49"
$ws.Range("E17").Value = "  -1.85%  "

$ws.Range("E18").Value = "  +0.04%  "

$ws.Range("D19").Value = "'6.40"
$ws.Range("E19").Value = "  -3.64%  "

$ws.Range("D20").Value = "'10.34"
$ws.Range("E20").Value = "  +4.50%  "

$ws.Range("E21").Value = "  -3.00%  "

$ws.Range("D22").Value = "'433.82"
$ws.Range("E22").Value = "  -4.10%  "

$ws.Range("D23").Value = "'0.609"
$ws.Range("E23").Value = "  -4.78%  "

$ws.Range("D24").Value = "'79.55"
$ws.Range("E24").Value = "  +1.67%  "

$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.626.79"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("E26").Value = "  -0.02%  "

$ws.Range("E27").Value = "  -8.62%  "

$ws.Range("D28").Value = "'9.81"
$ws.Range("E28").Value = "  -4.22%  "

$ws.Range("D29").Value = "'8.22"
$ws.Range("E29").Value = "  -9.61%  "

$ws.Range("D30").Value = "'2.49"
$ws.Range("E30").Value = "  -1.62%  "

$ws.Range("D31").Value = "'1.60"
$ws.Range("E31").Value = "  -5.28%  "

$ws.Range("E32").Value = "  +0.15%  "

$ws.Range("D33").Value = "'0.165"
$ws.Range("E33").Value = "  -3.85%  "

$ws.Range("D34").Value = "'25.39"
$ws.Range("E34").Value = "  -2.20%  "

$ws.Range("D35").Value = "3.483.34"
$ws.Range("E35").Value = "  -1.38%  "

$ws.Range("E36").Value = "  -5.21%  "

$ws.Range("D37").Value = "'1.80"
$ws.Range("E37").Value = "  -5.09%  "

$ws.Range("D38").Value = "'7.98"
$ws.Range("E38").Value = "  -1.91%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("E40").Value = "  -0.06%  "

$ws.Range("D41").Value = "'0.0891"
$ws.Range("E41").Value = "  -2.09%  "

$ws.Range("D42").Value = "'169.96"
$ws.Range("E42").Value = "  -2.92%  "

$ws.Range("D43").Value = "'5.42"
$ws.Range("E43").Value = "  -2.83%  "

$ws.Range("E44").Value = "  -11.36%  "

$ws.Range("E45").Value = "  +0.94%  "

$ws.Range("D46").Value = "'28.88"
$ws.Range("E46").Value = "  -6.65%  "

$ws.Range("D47").Value = "'45.70"
$ws.Range("E47").Value = "  -1.48%  "

$ws.Range("D48").Value = "'1.23"
$ws.Range("E48").Value = "  -7.75%  "

$ws.Range("E49").Value = "  -3.53%  "

$ws.Range("D50").Value = "'2.41"
$ws.Range("E50").Value = "  -6.57%  "

$ws.Range("D51").Value = "'0.966"
$ws.Range("E51").Value = "  -3.78%  "
